$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.18"
$ws.Range("G2").Value = "'22"
$ws.Range("D3").Value = "'23.11"
$ws.Range("G3").Value = "'22"
$ws.Range("D4").Value = "'5.380"
$ws.Range("G4").Value = "'22"
$ws.Range("D5").Value = "'0.05978"
$ws.Range("G5").Value = "'22"
$ws.Range("G6").Value = "'22"
$ws.Range("D7").Value = "'6.493"
$ws.Range("G7").Value = "'22"
$ws.Range("D8").Value = "'0.8091"
$ws.Range("G8").Value = "'22"
$ws.Range("D9").Value = "'0.9070"
$ws.Range("G9").Value = "'22"
$ws.Range("D10").Value = "'0.1421"
$ws.Range("G10").Value = "'22"
$ws.Range("D11").Value = "'0.07413"
$ws.Range("G11").Value = "'22"
$ws.Range("D12").Value = "'0.03311"
$ws.Range("G12").Value = "'22"
$ws.Range("D13").Value = "'0.03070"
$ws.Range("G13").Value = "'22"
$ws.Range("D14").Value = "'0.09336"
$ws.Range("G14").Value = "'22"
$ws.Range("D15").Value = "'3.849"
$ws.Range("G15").Value = "'22"
$ws.Range("D16").Value = "'0.001596"
$ws.Range("G16").Value = "'22"
$ws.Range("D17").Value = "'0.04626"
$ws.Range("G17").Value = "'22"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005932"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").Value = "'22"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006102"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("G19").Value = "'22"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.005040"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("G20").Value = "'22"
$ws.Range("B21").Value = "UpBots"
$ws.Range("C21").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D21").Value = "'0.007494"
$ws.Range("E21").Value = "20UpBotsUBXTBestin24h"
$ws.Range("G21").Value = "'22"
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").Value = "'0.0009799"
$ws.Range("E22").Value = "21BitKanKAN"
$ws.Range("G22").Value = "'22"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").Value = "'0.00007793"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("G23").Value = "'22"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").Value = "'3.616"
$ws.Range("E24").Value = "23LEOLEO"
$ws.Range("G24").Value = "'22"
$ws.Range("B25").Value = "BTSEToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D25").Value = "'2.163"
$ws.Range("E25").Value = "24BTSETokenBTSE"
$ws.Range("G25").Value = "'22"
$ws.Range("D26").Value = "'0.3214"
$ws.Range("G26").Value = "'22"
$ws.Range("G27").Value = "'22"
$ws.Range("G28").Value = "'22"
$ws.Range("G29").Value = "'22"
$ws.Range("G30").Value = "'22"
$ws.Range("G31").Value = "'22"
$ws.Range("G32").Value = "'22"
$ws.Range("G33").Value = "'22"
$ws.Range("G34").Value = "'22"
$ws.Range("G35").Value = "'22"
$ws.Range("G36").Value = "'22"
$ws.Range("G37").Value = "'22"
$ws.Range("G38").Value = "'22"
$ws.Range("G39").Value = "'22"
$ws.Range("G40").Value = "'22"
$ws.Range("D41").Value = "'0.006190"
$ws.Range("G41").Value = "'22"
$ws.Range("G42").Value = "'22"
$ws.Range("D43").Value = "'0.002797"
$ws.Range("G43").Value = "'22"
$ws.Range("D44").Value = "'0.007164"
$ws.Range("G44").Value = "'22"
$ws.Range("D45").Value = "'0.00005180"
$ws.Range("G45").Value = "'22"
$ws.Range("G46").Value = "'22"
$ws.Range("D47").Value = "'0.0005792"
$ws.Range("G47").Value = "'22"
$ws.Range("D48").Value = "'1.044"
$ws.Range("G48").Value = "'22"
$ws.Range("D49").Value = "'0.002262"
$ws.Range("G49").Value = "'22"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("G50").Value = "'22"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("G51").Value = "'22"
